$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3060.25
$ws.Range("I40").Value = 1350
$ws.Range("J40").Value = 3630.3333
$ws.Range("K40").Value = 1350
$ws.Range("L40").Value = 3630.3333
$ws.Range("M40").Value = -1175
$ws.Range("N40").Value = -3980.3333
# Row 51
$ws.Range("H51").Value = 2654.2856
$ws.Range("I51").Value = 2600
$ws.Range("J51").Value = 2663.3333
$ws.Range("K51").Value = 2600
$ws.Range("L51").Value = 2663.3333
$ws.Range("M51").Value = -2116
$ws.Range("N51").Value = -3631.3333
# Row 94
$ws.Range("H94").Value = 3900
$ws.Range("I94").Value = 3900
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3900
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -3449
$ws.Range("N94").Value = ""
# Row 98
$ws.Range("H98").Value = 1480.625
$ws.Range("I98").Value = 1262.2727
$ws.Range("J98").Value = 1961
$ws.Range("K98").Value = 1262.2727
$ws.Range("L98").Value = 1961
$ws.Range("M98").Value = 235.7273
$ws.Range("N98").Value = -4957
# Row 116
$ws.Range("H116").Value = 2961007.5
$ws.Range("I116").Value = 10991912
$ws.Range("J116").Value = 2253.4736
$ws.Range("K116").Value = 10991912
$ws.Range("L116").Value = 2253.4736
$ws.Range("M116").Value = -10988470
$ws.Range("N116").Value = -9137.473599999999
# Row 122
$ws.Range("H122").Value = 1480.625
$ws.Range("I122").Value = 1262.2727
$ws.Range("J122").Value = 1961
$ws.Range("K122").Value = 3786.8181
$ws.Range("L122").Value = 5883
$ws.Range("M122").Value = -1336.8181
$ws.Range("N122").Value = -10783
# Row 129
$ws.Range("H129").Value = 714.7449
$ws.Range("I129").Value = 606.5
$ws.Range("K129").Value = 1819.5
$ws.Range("M129").Value = 3180.5
# Row 132
$ws.Range("H132").Value = 3179.5
$ws.Range("I132").Value = 2520.4883
$ws.Range("J132").Value = 5755.636
$ws.Range("K132").Value = 7561.4649
$ws.Range("L132").Value = 17266.908
$ws.Range("M132").Value = -5031.4649
$ws.Range("N132").Value = -22326.908
# Row 138
$ws.Range("H138").Value = 2161.1958
$ws.Range("I138").Value = 938
$ws.Range("J138").Value = 2478.9092
$ws.Range("K138").Value = 2814
$ws.Range("L138").Value = 7436.7276
$ws.Range("M138").Value = 2326
$ws.Range("N138").Value = -17716.7276

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9264.690000000001
$ws.Range("I32").Value = 5644.7896
$ws.Range("J32").Value = 20727.709
$ws.Range("K32").Value = 5644.7896
$ws.Range("L32").Value = 20727.709
$ws.Range("M32").Value = -5357.7896
$ws.Range("N32").Value = -21301.709
# Row 61
$ws.Range("H61").Value = 3057.5881
$ws.Range("I61").Value = 2505.8462
$ws.Range("J61").Value = 4850.75
$ws.Range("K61").Value = 2505.8462
$ws.Range("L61").Value = 4850.75
$ws.Range("M61").Value = -2293.8462
$ws.Range("N61").Value = -5274.75
# Row 132
$ws.Range("H132").Value = 1186.4286
$ws.Range("I132").Value = 648.0513
$ws.Range("J132").Value = 2421.5293
$ws.Range("K132").Value = 1944.1539
$ws.Range("L132").Value = 7264.5879
$ws.Range("M132").Value = 585.8461000000002
$ws.Range("N132").Value = -12324.5879
# Row 136
$ws.Range("H136").Value = 3057.5881
$ws.Range("I136").Value = 2505.8462
$ws.Range("J136").Value = 4850.75
$ws.Range("K136").Value = 7517.5386
$ws.Range("L136").Value = 14552.25
$ws.Range("M136").Value = -4967.5386
$ws.Range("N136").Value = -19652.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2930.2046
$ws.Range("I31").Value = 2010.6522
$ws.Range("J31").Value = 3937.3333
$ws.Range("K31").Value = 2010.6522
$ws.Range("L31").Value = 3937.3333
$ws.Range("M31").Value = -1715.6522
$ws.Range("N31").Value = -4527.3333
# Row 34
$ws.Range("H34").Value = 2930.2046
$ws.Range("I34").Value = 2010.6522
$ws.Range("J34").Value = 3937.3333
$ws.Range("K34").Value = 2010.6522
$ws.Range("L34").Value = 3937.3333
$ws.Range("M34").Value = -1808.6522
$ws.Range("N34").Value = -4341.3333
# Row 58
$ws.Range("H58").Value = 2117.25
$ws.Range("I58").Value = 1875.9412
$ws.Range("J58").Value = 2295.6086
$ws.Range("K58").Value = 1875.9412
$ws.Range("L58").Value = 2295.6086
$ws.Range("M58").Value = -1672.9412
$ws.Range("N58").Value = -2701.6086
# Row 134
$ws.Range("H134").Value = 2890.0588
$ws.Range("I134").Value = 1599.6818
$ws.Range("J134").Value = 5255.75
$ws.Range("K134").Value = 4799.0454
$ws.Range("L134").Value = 15767.25
$ws.Range("M134").Value = -2264.0454
$ws.Range("N134").Value = -20837.25
# Row 136
$ws.Range("H136").Value = 2117.25
$ws.Range("I136").Value = 1875.9412
$ws.Range("J136").Value = 2295.6086
$ws.Range("K136").Value = 5627.8236
$ws.Range("L136").Value = 6886.825800000001
$ws.Range("M136").Value = -3077.8236
$ws.Range("N136").Value = -11986.8258

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1567861.9
$ws.Range("I113").Value = 6896912
$ws.Range("J113").Value = 494.11765
$ws.Range("K113").Value = 20690736
$ws.Range("L113").Value = 1482.35295
$ws.Range("M113").Value = -20688566
$ws.Range("N113").Value = -5822.35295
# Row 122
$ws.Range("H122").Value = 788021.4399999999
$ws.Range("I122").Value = 5166.2607
$ws.Range("J122").Value = 3788966.5
$ws.Range("K122").Value = 46496.3463
$ws.Range("L122").Value = 34100698.5
$ws.Range("M122").Value = -44046.3463
$ws.Range("N122").Value = -34105598.5
# Row 131
$ws.Range("H131").Value = 775.3514
$ws.Range("I131").Value = 411.1
$ws.Range("J131").Value = 910.2593000000001
$ws.Range("K131").Value = 1233.3
$ws.Range("L131").Value = 2730.7779
$ws.Range("M131").Value = 3806.7
$ws.Range("N131").Value = -12810.7779
# Row 132
$ws.Range("H132").Value = 11852999
$ws.Range("I132").Value = 1350
$ws.Range("J132").Value = 13676329
$ws.Range("K132").Value = 12150
$ws.Range("L132").Value = 123086961
$ws.Range("M132").Value = -9620
$ws.Range("N132").Value = -123092021

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2398.6597
$ws.Range("I132").Value = 1956.6471
$ws.Range("J132").Value = 3554.6924
$ws.Range("K132").Value = 5869.9413
$ws.Range("L132").Value = 10664.0772
$ws.Range("M132").Value = -3339.9413
$ws.Range("N132").Value = -15724.0772

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2162.82
$ws.Range("I132").Value = 1945.0834
$ws.Range("J132").Value = 3305.9375
$ws.Range("K132").Value = 5835.2502
$ws.Range("L132").Value = 9917.8125
$ws.Range("M132").Value = -3305.2502
$ws.Range("N132").Value = -14977.8125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1831.4147
$ws.Range("I126").Value = 1979.32
$ws.Range("J126").Value = 1600.3125
$ws.Range("K126").Value = 5937.96
$ws.Range("L126").Value = 4800.9375
$ws.Range("M126").Value = -3467.96
$ws.Range("N126").Value = -9740.9375

